$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# Update header row text (row 1) to new wording / capitalization
$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Cluster number"
$ws.Range("C1").Value = "Average family income (`$2012), 2000"
$ws.Range("D1").Value = "Average family income (`$2012), 2008–12"
$ws.Range("E1").Value = "Change in family income, 2000 to 2008–12"
$ws.Range("F1").Value = "Homeownership rate, 2008–12"
$ws.Range("G1").Value = "Percent below poverty level, 2000"
$ws.Range("H1").Value = "Percent below poverty, 2008–12"
$ws.Range("I1").Value = "Unemployment rate, 2000"
$ws.Range("J1").Value = "Unemployment rate, 2008–12"
$ws.Range("K1").Value = "Median sales price of single family home (`$2012), 2000"
$ws.Range("L1").Value = "Median sales price of single family home (`$2012), 2012"

# Resize the workbook window (bookViews/workbookView)
$excel.Width = 31940
$excel.Height = 20260
